{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n// the \"\u00a9 2020 . Contact: ...\" footer line, and the blank paragraph that\n// separated them from the \"LOQ4205: Sistemas Produtivos II (Requisito\n// fraco)\" line, leaving that requirement line directly followed by the\n// single blank paragraph that used to sit just before the page break.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"LOQ4205: Sistemas Produtivos II (Requisito fraco)\") {\n    anchorIndex = i;\n  }\n}\n\nif (anchorIndex !== -1) {\n  const blank = items[anchorIndex + 1];\n  const jupiter = items[anchorIndex + 2];\n  const footer = items[anchorIndex + 3];\n\n  if (footer && footer.text.indexOf(\"luizeleno@usp.br\") !== -1) {\n    footer.delete();\n  }\n  if (jupiter && jupiter.text === \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n    jupiter.delete();\n  }\n  if (blank && blank.text === \"\") {\n    blank.delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that holds the \"Requisitos\" section's course-code\n# line, then delete the three trailing paragraphs that followed it:\n#   1) a blank paragraph\n#   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3) the \"(c) 2020 . Contact: ...\" footer line\n# We search by text so the script is resilient to any pre-existing\n# paragraph-count differences, deleting from the bottom up so indices\n# stay valid as each paragraph disappears.\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\") -eq \"LOQ4205: Sistemas Produtivos II (Requisito fraco)\") {\n        $target = $p\n    }\n}\n\nif ($target -ne $null) {\n    $idx = $target.Index\n\n    $p3 = $d.Paragraphs($idx + 3)\n    if ($p3.Range.Text.TrimEnd(\"`r\") -like \"*2020*luizeleno@usp.br*\") {\n        $p3.Range.Delete()\n    }\n\n    $p2 = $d.Paragraphs($idx + 2)\n    if ($p2.Range.Text.TrimEnd(\"`r\") -eq \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n        $p2.Range.Delete()\n    }\n\n    $p1 = $d.Paragraphs($idx + 1)\n    if ($p1.Range.Text.TrimEnd(\"`r\") -eq \"\") {\n        $p1.Range.Delete()\n    }\n}\n"}
